$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helper: split the run containing [splitStart, splitEnd) away from its
# neighbours by toggling a character-formatting property on and back off.
# This forces the host's run-recompute step to keep the touched span as
# its own <w:r> instead of re-coalescing it with identically formatted
# neighbours.
# ---------------------------------------------------------------------
function Split-Run([int]$splitStart, [int]$splitEnd) {
    $r = $d.Range($splitStart, $splitEnd)
    $r.Font.Bold = $true
    $r.Font.Bold = $false
}

# =======================================================================
# 1) "(c)Paul Scharfenberger 2004-2016"  ->
#    "(c)" / "2004-2016 " / "Paul Scharfenberger"   (three runs)
# =======================================================================
$copyright = [char]0x00A9
$oldLine1 = $copyright + "Paul Scharfenberger 2004-2016"
$newLine1 = $copyright + "2004-2016 Paul Scharfenberger"

$found = $d.Content.Find.Execute($oldLine1, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Could not find copyright line" }
$line1Range = $d.Range($d.Content.Find.Parent.Start, $d.Content.Find.Parent.End)

# Re-find precisely (Find.Execute already collapsed Content to the match)
$m1 = $d.Content
$m1.Find.Execute($oldLine1, $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$p1Start = $m1.Start
$p1End = $m1.End

# Replace the whole paragraph's text (minus the trailing paragraph mark)
# with the rearranged text, in one run.
$textRange = $d.Range($p1Start, $p1End)
$textRange.Delete()
$insPoint1 = $d.Range($p1Start, $p1Start)
$insPoint1.InsertAfter($newLine1)

# Split points, measured from the paragraph start:
#   "(c)" = 1 char, "2004-2016 " = 10 chars, "Paul Scharfenberger" = 20 chars
$midStart = $p1Start + 1
$midEnd = $p1Start + 1 + 10
Split-Run $midStart $midEnd

# =======================================================================
# 2) "2 April 2016" -> "2 " (unchanged) / "May" / " 2016"
# =======================================================================
$oldDate = "April 2016"
$m2 = $d.Content
$found2 = $m2.Find.Execute($oldDate, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) { throw "Could not find date line" }
$p2Start = $m2.Start
$p2End = $m2.End

$dateRange = $d.Range($p2Start, $p2End)
$dateRange.Delete()
$insPoint2 = $d.Range($p2Start, $p2Start)
$insPoint2.InsertAfter("May 2016")

# Split point: after "May" (3 chars)
$maySplitStart = $p2Start
$maySplitEnd = $p2Start + 3
Split-Run $maySplitStart $maySplitEnd

Write-Output "done"
